$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.963.97'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '2.469.19'
$ws.Range("E3").Value = '  -2.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.09'
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.20'
$ws.Range("E6").Value = '  -1.85%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  -2.21%  '
$ws.Range("D9").Value = '2.469.77'
$ws.Range("E9").Value = '  -2.58%  '
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.165'
$ws.Range("E11").Value = '  -1.04%  '
$ws.Range("E12").Value = '  -2.28%  '
$ws.Range("D14").Value = '2.919.56'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.21'
$ws.Range("E15").Value = '  -3.85%  '
$ws.Range("D16").Value = '66.854.05'
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("E17").Value = '  -3.74%  '
$ws.Range("D18").Value = '2.452.31'
$ws.Range("E18").Value = '  -3.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.94'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.38'
$ws.Range("E20").Value = '  -9.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '349.57'
$ws.Range("E21").Value = '  -4.37%  '
$ws.Range("E22").Value = '  -3.80%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.42'
$ws.Range("E24").Value = '  -5.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.19'
$ws.Range("E26").Value = '  -3.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.08'
$ws.Range("E27").Value = '  -7.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -37.99%  '
$ws.Range("D29").Value = '2.585.07'
$ws.Range("E29").Value = '  -2.90%  '
$ws.Range("D30").Value = '0.0₃0896'
$ws.Range("E30").Value = '  -5.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '507.21'
$ws.Range("E31").Value = '  -5.73%  '
$ws.Range("E32").Value = '  -7.64%  '
$ws.Range("E33").Value = '  -5.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.22'
$ws.Range("E34").Value = '  -5.21%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '158.63'
$ws.Range("E36").Value = '  -0.67%  '
$ws.Range("E37").Value = '  -11.08%  '
$ws.Range("E39").Value = '  -5.80%  '
$ws.Range("E40").Value = '  -7.90%  '
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.67'
$ws.Range("E42").Value = '  -6.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.79'
$ws.Range("E43").Value = '  -5.47%  '
$ws.Range("E44").Value = '  -5.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.35'
$ws.Range("E45").Value = '  -4.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.81'
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.96'
$ws.Range("E47").Value = '  -5.31%  '
$ws.Range("E48").Value = '  -7.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.509'
$ws.Range("E49").Value = '  -7.16%  '
$ws.Range("D50").Value = '0.0₆0251'
$ws.Range("E50").Value = '  -8.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0728'
$ws.Range("E51").Value = '  -1.94%  '
